$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 19.36022366666667
$ws.Range("H2").Value = 58.080671
$ws.Range("I2").Value = 0.005884129141485179
$ws.Range("J2").Value = 0.005884129141485179
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 31.82741333333333
$ws.Range("N2").Value = 95.48223999999999
$ws.Range("O2").Value = 0.114390792932228
$ws.Range("P2").Value = 0.114390792932228
$ws.Range("Q2").Value = 616.1858408647821
$ws.Range("R2").Value = 5545.672567783039
$ws.Range("S2").Value = 0.0006730901982101196
$ws.Range("T2").Value = 0.0006730901982101198

# row3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 19.36022366666667
$ws.Range("H3").Value = 58.080671
$ws.Range("I3").Value = 0.005884129141485179
$ws.Range("J3").Value = 0.005884129141485179
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 85.46317833333335
$ws.Range("N3").Value = 256.389535
$ws.Range("O3").Value = 0.307162904935779
$ws.Range("P3").Value = 0.307162904935779
$ws.Range("Q3").Value = 1654.586247797554
$ws.Range("R3").Value = 14891.27623017799
$ws.Range("S3").Value = 0.001807386200115859
$ws.Range("T3").Value = 0.001807386200115859

# row4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 19.36022366666667
$ws.Range("H4").Value = 58.080671
$ws.Range("I4").Value = 0.005884129141485179
$ws.Range("J4").Value = 0.005884129141485179
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 122.2478306666667
$ws.Range("N4").Value = 366.743492
$ws.Range("O4").Value = 0.4393704929064738
$ws.Range("P4").Value = 0.4393704929064738
$ws.Range("Q4").Value = 2366.745344471459
$ws.Range("R4").Value = 21300.70810024313
$ws.Range("S4").Value = 0.002585312721219689
$ws.Range("T4").Value = 0.00258531272121969

# row5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 19.36022366666667
$ws.Range("H5").Value = 58.080671
$ws.Range("I5").Value = 0.005884129141485179
$ws.Range("J5").Value = 0.005884129141485179
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 38.69562533333333
$ws.Range("N5").Value = 116.086876
$ws.Range("O5").Value = 0.1390758092255191
$ws.Range("P5").Value = 0.1390758092255191
$ws.Range("Q5").Value = 749.1559613748661
$ws.Range("R5").Value = 6742.403652373795
$ws.Range("S5").Value = 0.0008183400219395103
$ws.Range("T5").Value = 0.0008183400219395103

# row6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 3161.845459
$ws.Range("H6").Value = 9485.536377
$ws.Range("I6").Value = 0.9609758299542277
$ws.Range("J6").Value = 0.9609758299542278
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 31.82741333333333
$ws.Range("N6").Value = 95.48223999999999
$ws.Range("O6").Value = 0.114390792932228
$ws.Range("P6").Value = 0.114390792932228
$ws.Range("Q6").Value = 100633.362319716
$ws.Range("R6").Value = 905700.2608774444
$ws.Range("S6").Value = 0.10992678717717
$ws.Range("T6").Value = 0.1099267871771701

# row7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 3161.845459
$ws.Range("H7").Value = 9485.536377
$ws.Range("I7").Value = 0.9609758299542277
$ws.Range("J7").Value = 0.9609758299542278
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 85.46317833333335
$ws.Range("N7").Value = 256.389535
$ws.Range("O7").Value = 0.307162904935779
$ws.Range("P7").Value = 0.307162904935779
$ws.Range("Q7").Value = 270221.3623249572
$ws.Range("R7").Value = 2431992.260924615
$ws.Range("S7").Value = 0.2951761275018118
$ws.Range("T7").Value = 0.2951761275018118

# row8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 3161.845459
$ws.Range("H8").Value = 9485.536377
$ws.Range("I8").Value = 0.9609758299542277
$ws.Range("J8").Value = 0.9609758299542278
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 122.2478306666667
$ws.Range("N8").Value = 366.743492
$ws.Range("O8").Value = 0.4393704929064738
$ws.Range("P8").Value = 0.4393704929064738
$ws.Range("Q8").Value = 386528.748266001
$ws.Range("R8").Value = 3478758.734394009
$ws.Range("S8").Value = 0.4222244240781968
$ws.Range("T8").Value = 0.4222244240781969

# row9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 3161.845459
$ws.Range("H9").Value = 9485.536377
$ws.Range("I9").Value = 0.9609758299542277
$ws.Range("J9").Value = 0.9609758299542278
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 38.69562533333333
$ws.Range("N9").Value = 116.086876
$ws.Range("O9").Value = 0.1390758092255191
$ws.Range("P9").Value = 0.1390758092255191
$ws.Range("Q9").Value = 122349.5872433654
$ws.Range("R9").Value = 1101146.285190288
$ws.Range("S9").Value = 0.1336484911970491
$ws.Range("T9").Value = 0.1336484911970491

# row10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 2.055785333333333
$ws.Range("H10").Value = 6.167356
$ws.Range("I10").Value = 0.0006248123263850286
$ws.Range("J10").Value = 0.0006248123263850286
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 31.82741333333333
$ws.Range("N10").Value = 95.48223999999999
$ws.Range("O10").Value = 0.114390792932228
$ws.Range("P10").Value = 0.114390792932228
$ws.Range("Q10").Value = 65.43032952860443
$ws.Range("R10").Value = 588.8729657574399
$ws.Range("S10").Value = 0.00007147277744901346
$ws.Range("T10").Value = 0.00007147277744901349

# row11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 2.055785333333333
$ws.Range("H11").Value = 6.167356
$ws.Range("I11").Value = 0.0006248123263850286
$ws.Range("J11").Value = 0.0006248123263850286
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 85.46317833333335
$ws.Range("N11").Value = 256.389535
$ws.Range("O11").Value = 0.307162904935779
$ws.Range("P11").Value = 0.307162904935779
$ws.Range("Q11").Value = 175.6939485577178
$ws.Range("R11").Value = 1581.24553701946
$ws.Range("S11").Value = 0.0001919191692121075
$ws.Range("T11").Value = 0.0001919191692121075

# row12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 2.055785333333333
$ws.Range("H12").Value = 6.167356
$ws.Range("I12").Value = 0.0006248123263850286
$ws.Range("J12").Value = 0.0006248123263850286
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 122.2478306666667
$ws.Range("N12").Value = 366.743492
$ws.Range("O12").Value = 0.4393704929064738
$ws.Range("P12").Value = 0.4393704929064738
$ws.Range("Q12").Value = 251.3152973163502
$ws.Range("R12").Value = 2261.837675847152
$ws.Range("S12").Value = 0.0002745240998178306
$ws.Range("T12").Value = 0.0002745240998178306

# row13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 2.055785333333333
$ws.Range("H13").Value = 6.167356
$ws.Range("I13").Value = 0.0006248123263850286
$ws.Range("J13").Value = 0.0006248123263850286
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 38.69562533333333
$ws.Range("N13").Value = 116.086876
$ws.Range("O13").Value = 0.1390758092255191
$ws.Range("P13").Value = 0.1390758092255191
$ws.Range("Q13").Value = 79.54989902442844
$ws.Range("R13").Value = 715.949091219856
$ws.Range("S13").Value = 0.00008689627990607702
$ws.Range("T13").Value = 0.00008689627990607702

# row14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 106.9830526666667
$ws.Range("H14").Value = 320.949158
$ws.Range("I14").Value = 0.03251522857790212
$ws.Range("J14").Value = 0.03251522857790212
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 31.82741333333333
$ws.Range("N14").Value = 95.48223999999999
$ws.Range("O14").Value = 0.114390792932228
$ws.Range("P14").Value = 0.114390792932228
$ws.Range("Q14").Value = 3404.993836883768
$ws.Range("R14").Value = 30644.94453195392
$ws.Range("S14").Value = 0.003719442779398864
$ws.Range("T14").Value = 0.003719442779398865

# row15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 106.9830526666667
$ws.Range("H15").Value = 320.949158
$ws.Range("I15").Value = 0.03251522857790212
$ws.Range("J15").Value = 0.03251522857790212
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 85.46317833333335
$ws.Range("N15").Value = 256.389535
$ws.Range("O15").Value = 0.307162904935779
$ws.Range("P15").Value = 0.307162904935779
$ws.Range("Q15").Value = 9143.111708695727
$ws.Range("R15").Value = 82288.00537826154
$ws.Range("S15").Value = 0.009987472064639276
$ws.Range("T15").Value = 0.009987472064639276

# row16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 106.9830526666667
$ws.Range("H16").Value = 320.949158
$ws.Range("I16").Value = 0.03251522857790212
$ws.Range("J16").Value = 0.03251522857790212
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 122.2478306666667
$ws.Range("N16").Value = 366.743492
$ws.Range("O16").Value = 0.4393704929064738
$ws.Range("P16").Value = 0.4393704929064738
$ws.Range("Q16").Value = 13078.44610659775
$ws.Range("R16").Value = 117706.0149593797
$ws.Range("S16").Value = 0.01428623200723952
$ws.Range("T16").Value = 0.01428623200723952

# row17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 106.9830526666667
$ws.Range("H17").Value = 320.949158
$ws.Range("I17").Value = 0.03251522857790212
$ws.Range("J17").Value = 0.03251522857790212
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 38.69562533333333
$ws.Range("N17").Value = 116.086876
$ws.Range("O17").Value = 0.1390758092255191
$ws.Range("P17").Value = 0.1390758092255191
$ws.Range("Q17").Value = 4139.7761230056
$ws.Range("R17").Value = 37257.98510705041
$ws.Range("S17").Value = 0.004522081726624463
$ws.Range("T17").Value = 0.004522081726624463
